$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2026-01-27 02:07:56"

for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}
